# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 7440
    3  = 7451
    8  = 118
    10 = 134
    13 = 675
    14 = 569
    19 = 79
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
